# Updated cryptos list on Sat Nov 30 07:35:14 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") "96.772.38"
Set-TextValue $ws.Range("E2") "  +0.79%  "
Set-TextValue $ws.Range("D3") "3.698.93"
Set-TextValue $ws.Range("E3") "  +4.26%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "244.06"
Set-TextValue $ws.Range("E5") "  +1.78%  "
Set-TextValue $ws.Range("D6") "1.93"
Set-TextValue $ws.Range("E6") "  +18.91%  "
Set-TextValue $ws.Range("D7") "675.09"
Set-TextValue $ws.Range("E7") "  +3.75%  "
Set-TextValue $ws.Range("D8") "0.426"
Set-TextValue $ws.Range("E8") "  +5.28%  "
Set-TextValue $ws.Range("E9") "  +5.89%  "
Set-TextValue $ws.Range("D11") "3.694.28"
Set-TextValue $ws.Range("E11") "  +4.16%  "
Set-TextValue $ws.Range("D12") "45.42"
Set-TextValue $ws.Range("E12") "  +4.87%  "
Set-TextValue $ws.Range("D13") "0.205"
Set-TextValue $ws.Range("E13") "  +1.90%  "
Set-TextValue $ws.Range("D14") "6.59"
Set-TextValue $ws.Range("E14") "  +3.69%  "
Set-TextValue $ws.Range("D15") "4.389.22"
Set-TextValue $ws.Range("E15") "  +4.30%  "
Set-TextValue $ws.Range("D16") "96.493.63"
Set-TextValue $ws.Range("E16") "  +0.59%  "
Set-TextValue $ws.Range("E17") "  +2.10%  "
Set-TextValue $ws.Range("E18") "  +13.70%  "
Set-TextValue $ws.Range("D19") "3.687.21"
Set-TextValue $ws.Range("E19") "  +3.78%  "
Set-TextValue $ws.Range("D20") "13.08"
Set-TextValue $ws.Range("E20") "  +5.55%  "
Set-TextValue $ws.Range("D21") "18.60"
Set-TextValue $ws.Range("E21") "  +5.54%  "
Set-TextValue $ws.Range("D22") "0.554"
Set-TextValue $ws.Range("E22") "  +6.11%  "
Set-TextValue $ws.Range("D23") "517.24"
Set-TextValue $ws.Range("E23") "  +2.57%  "
Set-TextValue $ws.Range("D24") "3.43"
Set-TextValue $ws.Range("E24") "  +1.42%  "
Set-TextValue $ws.Range("D25") "0.0000211"
Set-TextValue $ws.Range("E25") "  +7.09%  "
Set-TextValue $ws.Range("E26") "  +0.88%  "
Set-TextValue $ws.Range("D27") "101.89"
Set-TextValue $ws.Range("E27") "  +6.50%  "
Set-TextValue $ws.Range("D28") "13.01"
Set-TextValue $ws.Range("E28") "  +2.91%  "
Set-TextValue $ws.Range("D29") "0.169"
Set-TextValue $ws.Range("E29") "  +11.96%  "
Set-TextValue $ws.Range("D30") "3.10"
Set-TextValue $ws.Range("E30") "  +4.34%  "
Set-TextValue $ws.Range("D31") "12.16"
Set-TextValue $ws.Range("E31") "  +7.46%  "
Set-TextValue $ws.Range("E32") "  -0.17%  "
Set-TextValue $ws.Range("E33") "  +2.82%  "
Set-TextValue $ws.Range("D34") "33.28"
Set-TextValue $ws.Range("E34") "  +6.77%  "
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  +0.13%  "
Set-TextValue $ws.Range("E36") "  +9.40%  "
Set-TextValue $ws.Range("D37") "0.595"
Set-TextValue $ws.Range("E37") "  +6.14%  "
Set-TextValue $ws.Range("D38") "8.82"
Set-TextValue $ws.Range("E38") "  +1.82%  "
Set-TextValue $ws.Range("D39") "617.08"
Set-TextValue $ws.Range("E39") "  +1.21%  "
Set-TextValue $ws.Range("D40") "42.32"
Set-TextValue $ws.Range("E40") "  +26.37%  "
Set-TextValue $ws.Range("D41") "0.161"
Set-TextValue $ws.Range("E41") "  +7.87%  "
Set-TextValue $ws.Range("D42") "0.966"
Set-TextValue $ws.Range("E42") "  +7.89%  "
Set-TextValue $ws.Range("E43") "  +8.71%  "
Set-TextValue $ws.Range("E44") "  -0.02%  "
Set-TextValue $ws.Range("E45") "  +7.89%  "
Set-TextValue $ws.Range("D46") "0.0449"
Set-TextValue $ws.Range("E46") "  +7.12%  "
Set-TextValue $ws.Range("E47") "  +25.55%  "
Set-TextValue $ws.Range("E48") "  +1.94%  "
Set-TextValue $ws.Range("E49") "  +0.34%  "
Set-TextValue $ws.Range("D50") "8.60"
Set-TextValue $ws.Range("D51") "54.61"
Set-TextValue $ws.Range("E51") "  +3.52%  "
